$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Return_with_prediction (G), return_pct_change (H) and mean_return_pct_change (I, only row 2)
# values for rows 2-57 to match the recomputed auto-recurrence results.

$ws.Range("G2").Value = 0.08694235855969235
$ws.Range("H2").Value = -8.212049018547841
$ws.Range("I2").Value = 20.6567289200182

$ws.Range("G3").Value = 0.07233998577375519
$ws.Range("H3").Value = 9.450310616614752

$ws.Range("G4").Value = 0.03999289774323695
$ws.Range("H4").Value = 326.5787348362544

$ws.Range("G5").Value = 0.005111523360389024
$ws.Range("H5").Value = 138.1747199289299

$ws.Range("G6").Value = -0.2276319691856447
$ws.Range("H6").Value = -2.909767923649343

$ws.Range("G7").Value = -0.1900678173357139
$ws.Range("H7").Value = 23.93505184428383

$ws.Range("G8").Value = -0.2819368714087576
$ws.Range("H8").Value = 23.85350390275953

$ws.Range("G9").Value = -0.3941028475273671
$ws.Range("H9").Value = 1.13676152601783

$ws.Range("G10").Value = -0.05174645253037699
$ws.Range("H10").Value = -419.3261914362334

$ws.Range("G11").Value = 0.04666627688234877
$ws.Range("H11").Value = 390.468456658705

$ws.Range("G12").Value = 0.2399401539179316
$ws.Range("H12").Value = 5.612667873161911

$ws.Range("G13").Value = 0.2395176471751634
$ws.Range("H13").Value = -9.046436121686098

$ws.Range("G14").Value = -0.01098870542907599
$ws.Range("H14").Value = -14.85724436502021

$ws.Range("G15").Value = -0.0007873914144341191
$ws.Range("H15").Value = -103.900542206657

$ws.Range("G16").Value = 0.1292567780830718
$ws.Range("H16").Value = 9.540472938861818

$ws.Range("G17").Value = 0.2280166062613969
$ws.Range("H17").Value = 4.196664621472414

$ws.Range("G18").Value = 0.04634951415079622
$ws.Range("H18").Value = -23.34691964877982

$ws.Range("G19").Value = 0.07526505416167774
$ws.Range("H19").Value = -16.45506100086438

$ws.Range("G20").Value = -0.1326812428403808
$ws.Range("H20").Value = 8.829551305471011

$ws.Range("G21").Value = -0.1787817762352114
$ws.Range("H21").Value = 10.53301743391644

$ws.Range("G22").Value = 0.05175422858573275
$ws.Range("H22").Value = -4.842984197706046

$ws.Range("G23").Value = 0.04933480619631943
$ws.Range("H23").Value = 20.79978159852176

$ws.Range("G24").Value = 0.1287733868815848
$ws.Range("H24").Value = 11.26637585906437

$ws.Range("G25").Value = 0.1511457619266983
$ws.Range("H25").Value = -0.6117115601626495

$ws.Range("G26").Value = 0.02983327101816859
$ws.Range("H26").Value = -43.57964461809772

$ws.Range("G27").Value = 0.05327939712502726
$ws.Range("H27").Value = 5.573147545398719

$ws.Range("G28").Value = 0.1565250412759756
$ws.Range("H28").Value = 2.363444817541757

$ws.Range("G29").Value = 0.1886110414476004
$ws.Range("H29").Value = 10.49070273401846

$ws.Range("G30").Value = 0.006772946851218597
$ws.Range("H30").Value = -65.38425442836699

$ws.Range("G31").Value = 0.0293334441588024
$ws.Range("H31").Value = 202.251756972065

$ws.Range("G32").Value = 0.02168183125740182
$ws.Range("H32").Value = -41.86127973940425

$ws.Range("G33").Value = 0.01112141651027694
$ws.Range("H33").Value = -57.39514875655669

$ws.Range("G34").Value = 0.09955773566323552
$ws.Range("H34").Value = -22.20107452649377

$ws.Range("G35").Value = 0.1514573519863638
$ws.Range("H35").Value = 17.71852468065908

$ws.Range("G36").Value = -0.006401564910383272
$ws.Range("H36").Value = -142.5865675604533

$ws.Range("G37").Value = 0.005124345175291234
$ws.Range("H37").Value = -66.53923369484622

$ws.Range("G38").Value = -0.03547773085368686
$ws.Range("H38").Value = -1637.049295859427

$ws.Range("G39").Value = -0.02552474354254478
$ws.Range("H39").Value = 23.59968490497057

$ws.Range("G40").Value = 0.1546745746402088
$ws.Range("H40").Value = 4.829299386763253

$ws.Range("G41").Value = 0.1312344590601114
$ws.Range("H41").Value = -18.68952628376135

$ws.Range("G42").Value = 0.05454724696877319
$ws.Range("H42").Value = -15.51563864802451

$ws.Range("G43").Value = 0.05109018977713883
$ws.Range("H43").Value = 46.97762025356295

$ws.Range("G44").Value = 0.01435164219241896
$ws.Range("H44").Value = 1.692504793060801

$ws.Range("G45").Value = 0.02124831308924308
$ws.Range("H45").Value = -48.2481046981845

$ws.Range("G46").Value = -0.02827203165444766
$ws.Range("H46").Value = 57.0463535916281

$ws.Range("G47").Value = -0.01197462943645966
$ws.Range("H47").Value = 71.01255002632715

$ws.Range("G48").Value = -0.1265432889806382
$ws.Range("H48").Value = -0.4507860386157451

$ws.Range("G49").Value = -0.1584548893929209
$ws.Range("H49").Value = 19.76212235255278

$ws.Range("G50").Value = 0.1301788591129626
$ws.Range("H50").Value = 19.56874943875523

$ws.Range("G51").Value = 0.1404532754250782
$ws.Range("H51").Value = 40.07391557574497

$ws.Range("G52").Value = 0.06976199292737169
$ws.Range("H52").Value = 17.01291634539596

$ws.Range("G53").Value = 0.06682192648760603
$ws.Range("H53").Value = -1.06710890358727

$ws.Range("G54").Value = -0.09505533539415029
$ws.Range("H54").Value = -35.94757787444652

$ws.Range("G55").Value = -0.1147419412313216
$ws.Range("H55").Value = -48.55979114590202

$ws.Range("G56").Value = 0.0563113151212511
$ws.Range("H56").Value = 22.88051474473359

$ws.Range("G57").Value = 0.1322624733189571
$ws.Range("H57").Value = 2458.13088066928
